$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38: Just Give Him a Serum
$ws.Range("H38").Value = 4393.952
$ws.Range("I38").Value = 2152.4666
$ws.Range("J38").Value = 9997.666999999999
$ws.Range("K38").Value = 6457.399800000001
$ws.Range("L38").Value = 29993.001
$ws.Range("M38").Value = -6085.399800000001
$ws.Range("N38").Value = -30737.001

# Row 138: All-night Crafting
$ws.Range("H138").Value = 20838294
$ws.Range("I138").Value = 1448.5714
$ws.Range("J138").Value = 29418172
$ws.Range("K138").Value = 4345.7142
$ws.Range("L138").Value = 88254516
$ws.Range("M138").Value = 794.2857999999997
$ws.Range("N138").Value = -88264796

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 1444.7826
$ws.Range("I141").Value = 1369.3636
$ws.Range("K141").Value = 4108.0908
$ws.Range("M141").Value = 1071.9092

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3181.4075
$ws.Range("I61").Value = 1870.5294
$ws.Range("J61").Value = 5409.9
$ws.Range("K61").Value = 1870.5294
$ws.Range("L61").Value = 5409.9
$ws.Range("M61").Value = -1658.5294
$ws.Range("N61").Value = -5833.9

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 53235.707
$ws.Range("I74").Value = 55462.74
$ws.Range("K74").Value = 55462.74
$ws.Range("M74").Value = -54588.74

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 53235.707
$ws.Range("I77").Value = 55462.74
$ws.Range("K77").Value = 277313.7
$ws.Range("M77").Value = -272945.7

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 1836.2051
$ws.Range("I122").Value = 1752.8438
$ws.Range("K122").Value = 5258.5314
$ws.Range("M122").Value = -2808.5314

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1971.46
$ws.Range("I132").Value = 2011.4615
$ws.Range("K132").Value = 6034.3845
$ws.Range("M132").Value = -3504.3845

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 3181.4075
$ws.Range("I136").Value = 1870.5294
$ws.Range("J136").Value = 5409.9
$ws.Range("K136").Value = 5611.5882
$ws.Range("L136").Value = 16229.7
$ws.Range("M136").Value = -3061.5882
$ws.Range("N136").Value = -21329.7

$ws = $wb.Worksheets.Item("BSM")
# Row 37: That's Some Fine Grinding
$ws.Range("H37").Value = 596.5714
$ws.Range("I37").Value = 785.2
$ws.Range("K37").Value = 785.2
$ws.Range("M37").Value = -648.2

# Row 96: Hammer Time
$ws.Range("H96").Value = 54235.848
$ws.Range("I96").Value = 6074.75
$ws.Range("K96").Value = 6074.75
$ws.Range("M96").Value = -3328.75

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 3389.8
$ws.Range("I99").Value = 1831.3334
$ws.Range("K99").Value = 1831.3334
$ws.Range("M99").Value = -333.3334

# Row 141: Awl Dreams Come True
$ws.Range("H141").Value = 42232.332
$ws.Range("J141").Value = 89988
$ws.Range("L141").Value = 89988
$ws.Range("N141").Value = -100348

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 5730938.5
$ws.Range("I58").Value = 1244.6666
$ws.Range("K58").Value = 1244.6666
$ws.Range("M58").Value = -1041.6666

# Row 131: An Integral Reward
$ws.Range("H131").Value = 29090.375
$ws.Range("J131").Value = 31531.857
$ws.Range("L131").Value = 31531.857
$ws.Range("N131").Value = -41611.857

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 43708.062
$ws.Range("I132").Value = 56651.35
$ws.Range("J132").Value = 3799.5833
$ws.Range("K132").Value = 169954.05
$ws.Range("L132").Value = 11398.7499
$ws.Range("M132").Value = -167424.05
$ws.Range("N132").Value = -16458.7499

# Row 136: Turali Quality
$ws.Range("H136").Value = 5730938.5
$ws.Range("I136").Value = 1244.6666
$ws.Range("K136").Value = 3733.9998
$ws.Range("M136").Value = -1183.9998

# Row 141: No Greater Treasure
$ws.Range("H141").Value = 83661.96000000001
$ws.Range("J141").Value = 90613.664
$ws.Range("L141").Value = 90613.664
$ws.Range("N141").Value = -100973.664

$ws = $wb.Worksheets.Item("CUL")
# Row 14: Keep Your Powder Dry
$ws.Range("H14").Value = 8528.416999999999
$ws.Range("I14").Value = 8528.416999999999
$ws.Range("K14").Value = 25585.251
$ws.Range("M14").Value = -25412.251

# Row 45: Don't Turn Up Your Nose
$ws.Range("H45").Value = 10000
$ws.Range("J45").Value = 10000
$ws.Range("L45").Value = 30000
$ws.Range("N45").Value = -31064

# Row 52: Made by Apple in Coerthas
$ws.Range("H52").Value = 4666.25
$ws.Range("J52").Value = 4666.25
$ws.Range("L52").Value = 13998.75
$ws.Range("N52").Value = -14530.75

# Row 108: Meet for Meat
$ws.Range("H108").Value = 1000
$ws.Range("I108").Value = 1000
$ws.Range("K108").Value = 3000
$ws.Range("M108").Value = -120

# Row 114: One Last Meal
$ws.Range("H114").Value = 7326.769
$ws.Range("I114").Value = 664.9
$ws.Range("K114").Value = 1994.7
$ws.Range("M114").Value = 1259.3

# Row 115: Mixology
$ws.Range("H115").Value = 5101.75
$ws.Range("I115").Value = 464.66666
$ws.Range("K115").Value = 1393.99998
$ws.Range("M115").Value = -218.9999800000001

# Row 128: A Historical Flavor
$ws.Range("H128").Value = 139990
$ws.Range("I128").Value = 139990
$ws.Range("K128").Value = 419970
$ws.Range("M128").Value = -414990

# Row 137: Creative Chocolate
$ws.Range("H137").Value = 8215
$ws.Range("J137").Value = 8306.799999999999
$ws.Range("L137").Value = 24920.4
$ws.Range("N137").Value = -35120.39999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 58: The Big Red
$ws.Range("H58").Value = 54000
$ws.Range("J58").Value = 54000
$ws.Range("L58").Value = 54000
$ws.Range("N58").Value = -54554

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 3751.9167
$ws.Range("I113").Value = 3345.6875
$ws.Range("K113").Value = 3345.6875
$ws.Range("M113").Value = -1175.6875

# Row 132: On Board for Lar
$ws.Range("H132").Value = 1059.1538
$ws.Range("I132").Value = 1067
$ws.Range("J132").Value = 1016
$ws.Range("K132").Value = 3201
$ws.Range("L132").Value = 3048
$ws.Range("M132").Value = -671
$ws.Range("N132").Value = -8108

$ws = $wb.Worksheets.Item("LTW")
# Row 56: Hold On Tight
$ws.Range("H56").Value = 41773
$ws.Range("I56").Value = 39051
$ws.Range("K56").Value = 39051
$ws.Range("M56").Value = -38360

# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 3765.0557
$ws.Range("J100").Value = 4249.25
$ws.Range("L100").Value = 4249.25
$ws.Range("N100").Value = -5331.25

# Row 122: Hell on Leather
$ws.Range("H122").Value = 3100.6365
$ws.Range("I122").Value = 2773.375
$ws.Range("J122").Value = 3973.3333
$ws.Range("K122").Value = 8320.125
$ws.Range("L122").Value = 11919.9999
$ws.Range("M122").Value = -5870.125
$ws.Range("N122").Value = -16819.9999

# Row 131: For What Was Gleaned
$ws.Range("H131").Value = 74998
$ws.Range("J131").Value = 74998
$ws.Range("L131").Value = 74998
$ws.Range("N131").Value = -85078

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 4734.07
$ws.Range("I132").Value = 2613.0908
$ws.Range("K132").Value = 7839.2724
$ws.Range("M132").Value = -5309.2724

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 1650.6904
$ws.Range("I136").Value = 1567
$ws.Range("K136").Value = 4701
$ws.Range("M136").Value = -2151

$ws = $wb.Worksheets.Item("WVR")
# Row 19: Dirt Cheap
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("N19").Value = 0

# Row 37: Bet You Anything
$ws.Range("H37").Value = 23495.666
$ws.Range("J37").Value = 23495.666
$ws.Range("L37").Value = 23495.666
$ws.Range("N37").Value = -23901.666

# Row 40: Our Man in Ul'dah
$ws.Range("H40").Value = 16183.143
$ws.Range("I40").Value = 7929
$ws.Range("J40").Value = 22373.75
$ws.Range("K40").Value = 7929
$ws.Range("L40").Value = 22373.75
$ws.Range("M40").Value = -7780
$ws.Range("N40").Value = -22671.75

# Row 61: Bundle Up, It's Odd out There
$ws.Range("H61").Value = 32873.75
$ws.Range("I61").Value = 29000
$ws.Range("K61").Value = 29000
$ws.Range("M61").Value = -28708

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 8089.5454
$ws.Range("I126").Value = 9141.714
$ws.Range("J126").Value = 6248.25
$ws.Range("K126").Value = 27425.142
$ws.Range("L126").Value = 18744.75
$ws.Range("M126").Value = -24955.142
$ws.Range("N126").Value = -23684.75

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 4547.3438
$ws.Range("I136").Value = 2821.9583
$ws.Range("K136").Value = 8465.874899999999
$ws.Range("M136").Value = -5915.874899999999

